# Generate Report for Handback
# - Marks the 33a0f370-... file's status as a failed handback transform
#   (instead of "Ready for handoff") on the Overview sheet and on each
#   per-locale sheet.
# - Records the handback/handoff filename-mismatch error detail for that
#   file on the zh-cn and de-de sheets.
# - Widens the "Error Detail" column (P) on those sheets so the message
#   is readable.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$newStatus = "Handback transform failed"

# Overview sheet: row 3 is the 33a0f370-... file; E = zh-cn status, F = de-de status
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

# Per-locale sheets: row 3 is the 33a0f370-... file; column C is Status
$wsZhCn.Range("C3").Value = $newStatus
$wsDeDe.Range("C3").Value = $newStatus

# Error Detail (column P) for the failed handback transform
$wsZhCn.Range("P3").Value = "Handback file name: t4ajssy4.lsk is different with handoff file name: 33a0f370-5b02-4d4a-afc0-b6f6c9d5f82b.caa85a7835fa7c8d902cda37d6f9b8520bfe31a3.zh-cn."
$wsDeDe.Range("P3").Value = "Handback file name: t4ajssy4.lsk is different with handoff file name: 33a0f370-5b02-4d4a-afc0-b6f6c9d5f82b.caa85a7835fa7c8d902cda37d6f9b8520bfe31a3.de-de."

# Widen the Error Detail column so the message is visible
$wsZhCn.Columns.Item(16).ColumnWidth = 39.17
$wsDeDe.Columns.Item(16).ColumnWidth = 39.17
